$d = $word.ActiveDocument

# Fix the wrong start date for the Research Fellow position: Jan -> Jun
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Jan 1998 " + [char]8211 + " Dec 1998", $true, $false, $false, $false, $false, $true, 1, $false, "Jun 1998 " + [char]8211 + " Dec 1998", 2)
